# "obstruent and stops fixes"
# Strip stray diacritics from a handful of phonetic transcriptions on the
# InputWords sheet (Yem and Abidji rows), and leave the selection where the
# author last left it (H6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InputWords")

# Yem row (row 4)
$ws.Range("B4").Value = "šuʔ-a"
$ws.Range("C4").Value = "toš-a"
$ws.Range("D4").Value = "kew-a"
$ws.Range("E4").Value = "kal-a"

# Nara row (row 5)
$ws.Range("C5").Value = "tıfini"

# Abidji row (row 6)
$ws.Range("B6").Value = "lu=bobu"
$ws.Range("C6").Value = "rʋwa"
$ws.Range("D6").Value = "n=di"
$ws.Range("E6").Value = "ahʋa"

$ws.Activate()
$ws.Range("H6").Select()
